$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was "Interstellar / 2018", now becomes "Look Back / 2023"
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Look Back"
$ws.Range("C2").Value = 2023
$ws.Range("D2").Value = "Action, Loneliness"
$ws.Range("G2").Value = "very sad, poor kid"

# Row 3: new row for "Interstellar / 2022" (id stays 2, was previously Taxi Driver row)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Interstellar"
$ws.Range("C3").Value = 2022
$ws.Range("D3").Value = "Fiction, Space"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "-"
$ws.Range("G3").Value = "-"

# Row 4: Taxi Driver moved down, year fixed to 1999, genre trimmed to Loneliness
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Taxi Driver"
$ws.Range("C4").Value = 1999
$ws.Range("D4").Value = "Loneliness"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = "-"

# Row 5: Oppenheimer
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Oppenheimer"
$ws.Range("C5").Value = 2023
$ws.Range("D5").Value = "History, Bomb, War"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "-"

# Row 6: The Batman
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "The Batman"
$ws.Range("C6").Value = 2022
$ws.Range("D6").Value = "action"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "-"
$ws.Range("G6").Value = "-"

# Row 7: Batman
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Batman"
$ws.Range("C7").Value = 2022
$ws.Range("D7").Value = "Action"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "-"
$ws.Range("G7").Value = "-"

# Row 8: Dawdaw
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Dawdaw"
$ws.Range("C8").Value = 234
$ws.Range("D8").Value = "Dawd, Ada, D, Wad"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "-"
$ws.Range("G8").Value = "-"

# Row 9: Sdawda
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Sdawda"
$ws.Range("C9").Value = 234
$ws.Range("D9").Value = "Awdawd, Awda, Ad, Aw"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "-"
$ws.Range("G9").Value = "-"
